$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, mirroring the style of D1 (bold/centered header style)
$ws.Range("E1").Value = "Colocação"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = $ws.Range("D1").HorizontalAlignment

# Ranking values for rows 2-8
$values = @("1º", "2º", "3º", "4º", "5º", "6º", "12º")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
